$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray legacy content that isn't part of the refined layout
$ws.Range("A5").ClearContents()

# ---------------------------------------------------------------------------
# Section 1: "System Function Points" mini-table (A1:B4)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "System Function Points"
$ws.Range("A1").Font.Bold = $true

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Must show exactly one item from the feed at a time"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Pagination device should conform to GEL"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Must accept feed in JSON format"

$ws.Range("A1:B4").BorderAround(1)

# ---------------------------------------------------------------------------
# Section 2: refined data-dictionary classification table (A6:C24)
# ---------------------------------------------------------------------------

# Header row
$ws.Range("A6").Value = "Tangible Things"
$ws.Range("B6").Value = "Related to system function point"
$ws.Range("C6").Value = "Tangible Things"

# Tangible things body
$ws.Range("A7").Value = "PaginationDevice"
$ws.Range("B7").Value = "1,2,3"
$ws.Range("C7").Value = "NewsTicker"

$ws.Range("A8").Value = "Page"
$ws.Range("B8").Value = "'1,"
$ws.Range("C8").Value = "NewsTickerItem?"

$ws.Range("A9").Value = "NewsTicker"
$ws.Range("B9").Value = "1, 2,3"
$ws.Range("C9").Value = "NewsTickerButton?"

$ws.Range("A10").Value = "NewsFeed"
$ws.Range("B10").Value = "'3,"
$ws.Range("C10").Value = "NewsFeed"

$ws.Range("A11").Value = "NewsFeedItem"
$ws.Range("B11").Value = "'1,"

$ws.Range("A12").Value = "PaginationDeviceButton"
$ws.Range("B12").Value = "2?"

# Roles section
$ws.Range("A14").Value = "Roles"
$ws.Range("A14").Font.Bold = $true

$ws.Range("A15").Value = "User?"
$ws.Range("B15").Value = "'0,"

# Events section
$ws.Range("A17").Value = "Events"
$ws.Range("A17").Font.Bold = $true
$ws.Range("C17").Value = "Events"
$ws.Range("C17").Font.Bold = $true

$ws.Range("A18").Value = "NewsFeedItemCycle"
$ws.Range("B18").Value = "'1,"
$ws.Range("C18").Value = "NewsFeedDidLoad"

$ws.Range("A19").Value = "NewsTickerWasCycled"
$ws.Range("B19").Value = "'1,"

$ws.Range("A20").Value = "NewsTickerLoadsFeed"
$ws.Range("B20").Value = "1,3"

# Interaction section
$ws.Range("A22").Value = "Interaction"
$ws.Range("A22").Font.Bold = $true

$ws.Range("A23").Value = "UserCyclesNewsTicker"
$ws.Range("B23").Value = "'1,"
$ws.Range("C23").Value = "NewsTickerBeginNavigation"

$ws.Range("A24").Value = "UserNavigatesNewsTicker"
$ws.Range("B24").Value = "'1,"
$ws.Range("A24").Font.Bold = $true
$ws.Range("B24").Font.Bold = $true
$ws.Range("C24").Font.Bold = $true

# blank but bordered cells
$ws.Range("C11").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("C15").Value = ""
$ws.Range("A16").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C19").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("A21").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("C24").Value = ""

# Borders: header row gets a full box on every cell
$ws.Range("A6:C6").Borders.LineStyle = 1

# Borders: body gets left/right on every cell (outer left/right edges plus
# an inside vertical line), and the final row also gets a bottom edge.
$ws.Range("A7:C24").Borders.Item(11).LineStyle = 1
$ws.Range("A7:C24").Borders.Item(7).LineStyle = 1
$ws.Range("A7:C24").Borders.Item(10).LineStyle = 1
$ws.Range("A24:C24").Borders.Item(9).LineStyle = 1

# Extra top border on C7 only
$ws.Range("C7").Borders.Item(8).LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths / selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.666666666666668
$ws.Columns.Item(2).ColumnWidth = 49.833333333333336
$ws.Columns.Item(3).ColumnWidth = 25.666666666666668

$null = $ws.Range("C10").Select()
